$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1100.0625
$ws.Range("J17").Value = 1166.7778
$ws.Range("L17").Value = 3500.3334
$ws.Range("N17").Value = -3836.3334
$ws.Range("H51").Value = 10259.406
$ws.Range("I51").Value = 10707.107
$ws.Range("J51").Value = 7125.5
$ws.Range("K51").Value = 10707.107
$ws.Range("L51").Value = 7125.5
$ws.Range("M51").Value = -10223.107
$ws.Range("N51").Value = -8093.5
$ws.Range("H80").Value = 2697.1667
$ws.Range("J80").Value = 3504.923
$ws.Range("L80").Value = 10514.769
$ws.Range("N80").Value = -12510.769
$ws.Range("H83").Value = 2697.1667
$ws.Range("J83").Value = 3504.923
$ws.Range("L83").Value = 31544.307
$ws.Range("N83").Value = -41528.307
$ws.Range("H107").Value = 813.86664
$ws.Range("I107").Value = 864.4286
$ws.Range("J107").Value = 106
$ws.Range("K107").Value = 864.4286
$ws.Range("L107").Value = 106
$ws.Range("M107").Value = -3946
$ws.Range("N107").Value = -3946

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 814.2
$ws.Range("I32").Value = 825.2273
$ws.Range("K32").Value = 825.2273
$ws.Range("M32").Value = -538.2273
$ws.Range("H45").Value = 1326.1578
$ws.Range("J45").Value = 1500
$ws.Range("L45").Value = 1500
$ws.Range("N45").Value = -2254
$ws.Range("H122").Value = 2169.7058
$ws.Range("I122").Value = 2169.7058
$ws.Range("K122").Value = 6509.117400000001
$ws.Range("M122").Value = -4059.117400000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4107.727
$ws.Range("I99").Value = 4340.5713
$ws.Range("J99").Value = 3700.25
$ws.Range("K99").Value = 4340.5713
$ws.Range("L99").Value = 3700.25
$ws.Range("M99").Value = -2842.5713
$ws.Range("N99").Value = -6696.25
$ws.Range("H107").Value = 1278.1428
$ws.Range("I107").Value = 987.25
$ws.Range("J107").Value = 1666
$ws.Range("K107").Value = 987.25
$ws.Range("L107").Value = 1666
$ws.Range("M107").Value = 932.75
$ws.Range("N107").Value = -5506

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1912
$ws.Range("I31").Value = 1165.7222
$ws.Range("J31").Value = 4598.6
$ws.Range("K31").Value = 1165.7222
$ws.Range("L31").Value = 4598.6
$ws.Range("M31").Value = -870.7221999999999
$ws.Range("N31").Value = -5188.6
$ws.Range("H34").Value = 1912
$ws.Range("I34").Value = 1165.7222
$ws.Range("J34").Value = 4598.6
$ws.Range("K34").Value = 1165.7222
$ws.Range("L34").Value = 4598.6
$ws.Range("M34").Value = -963.7221999999999
$ws.Range("N34").Value = -5002.6
$ws.Range("H99").Value = 3199.6667
$ws.Range("I99").Value = 2999.5
$ws.Range("K99").Value = 2999.5
$ws.Range("M99").Value = -1501.5
$ws.Range("H122").Value = 2499.6667
$ws.Range("I122").Value = 2499
$ws.Range("J122").Value = 2499.8
$ws.Range("K122").Value = 7497
$ws.Range("L122").Value = 7499.400000000001
$ws.Range("M122").Value = -5047
$ws.Range("N122").Value = -12399.4
$ws.Range("H126").Value = 3199.6667
$ws.Range("I126").Value = 2999.5
$ws.Range("K126").Value = 8998.5
$ws.Range("M126").Value = -6528.5
$ws.Range("H134").Value = 3849059.2
$ws.Range("I134").Value = 2641
$ws.Range("K134").Value = 7923
$ws.Range("M134").Value = -5388

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 112.5
$ws.Range("I2").Value = 88.14286
$ws.Range("J2").Value = 146.6
$ws.Range("K2").Value = 528.85716
$ws.Range("L2").Value = 879.5999999999999
$ws.Range("M2").Value = -415.85716
$ws.Range("N2").Value = -1105.6
$ws.Range("H4").Value = 137317.05
$ws.Range("I4").Value = 521.875
$ws.Range("J4").Value = 215485.72
$ws.Range("K4").Value = 1565.625
$ws.Range("L4").Value = 646457.16
$ws.Range("M4").Value = -1453.625
$ws.Range("N4").Value = -646681.16
$ws.Range("H8").Value = 819.5
$ws.Range("I8").Value = 819.5
$ws.Range("K8").Value = 2458.5
$ws.Range("M8").Value = -2319.5
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 464.3
$ws.Range("I107").Value = 441.58334
$ws.Range("K107").Value = 1324.75002
$ws.Range("M107").Value = 595.2499800000001
$ws.Range("H112").Value = 3327.9285
$ws.Range("I112").Value = 1058.5
$ws.Range("K112").Value = 3175.5
$ws.Range("M112").Value = -2067.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 9663.333000000001
$ws.Range("J42").Value = 9663.333000000001
$ws.Range("L42").Value = 9663.333000000001
$ws.Range("N42").Value = -10789.333
$ws.Range("H49").Value = 9663.333000000001
$ws.Range("J49").Value = 9663.333000000001
$ws.Range("L49").Value = 9663.333000000001
$ws.Range("N49").Value = -9957.333000000001
$ws.Range("H68").Value = 2171.2727
$ws.Range("I68").Value = 1990.7142
$ws.Range("J68").Value = 2487.25
$ws.Range("K68").Value = 1990.7142
$ws.Range("L68").Value = 2487.25
$ws.Range("M68").Value = -1241.7142
$ws.Range("N68").Value = -3985.25
$ws.Range("H71").Value = 2171.2727
$ws.Range("I71").Value = 1990.7142
$ws.Range("J71").Value = 2487.25
$ws.Range("K71").Value = 9953.571
$ws.Range("L71").Value = 12436.25
$ws.Range("M71").Value = -6209.571
$ws.Range("N71").Value = -19924.25
$ws.Range("H122").Value = 3386.9375
$ws.Range("J122").Value = 3585.4
$ws.Range("L122").Value = 10756.2
$ws.Range("N122").Value = -15656.2
$ws.Range("H132").Value = 4388.1577
$ws.Range("I132").Value = 4569.6924
$ws.Range("K132").Value = 13709.0772
$ws.Range("M132").Value = -11179.0772

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 50000
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51498
$ws.Range("H72").Value = 50000
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -157488
$ws.Range("H132").Value = 1975.6285
$ws.Range("I132").Value = 1975.6285
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5926.8855
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3396.8855
$ws.Range("N132").ClearContents()
